# This script applies the ITA model update of 2025-08-10 22:40.
# It reorders the cost-class rows within several solar/wind CF-class
# resource blocks on the "solar" and "wind" sheets by swapping the
# process name/description (columns C, D, K) and the associated
# numeric attributes (columns M, N, O, P - capacity, cap factor,
# investment cost, cost-class rank) between the affected rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "solar" (workbook sheet5.xml): spv-ITA_13 cost class 3 & 4
# ---------------------------------------------------------------
$wsSolar = $wb.Worksheets.Item("solar")

$wsSolar.Range("C30").Value = "e_spv-ITA_13_c4"
$wsSolar.Range("D30").Value = "solar resource -- CF class spv-ITA_13 -- cost class 4"
$wsSolar.Range("K30").Value = "e_spv-ITA_13_c4"
$wsSolar.Range("P30").Value = 4

$wsSolar.Range("C31").Value = "e_spv-ITA_13_c3"
$wsSolar.Range("D31").Value = "solar resource -- CF class spv-ITA_13 -- cost class 3"
$wsSolar.Range("K31").Value = "e_spv-ITA_13_c3"
$wsSolar.Range("P31").Value = 3

# ---------------------------------------------------------------
# Sheet "wind" (workbook sheet6.xml)
# ---------------------------------------------------------------
$wsWind = $wb.Worksheets.Item("wind")

# won-ITA_32 cost class 1 & 3 (row 7 <-> row 8)
$wsWind.Range("C7").Value = "e_won-ITA_32_c3"
$wsWind.Range("D7").Value = "wind resource -- CF class won-ITA_32 -- cost class 3"
$wsWind.Range("K7").Value = "e_won-ITA_32_c3"
$wsWind.Range("M7").Value = 2.4990000000000001
$wsWind.Range("O7").Value = 39.630069093581724
$wsWind.Range("P7").Value = 3

$wsWind.Range("C8").Value = "e_won-ITA_32_c1"
$wsWind.Range("D8").Value = "wind resource -- CF class won-ITA_32 -- cost class 1"
$wsWind.Range("K8").Value = "e_won-ITA_32_c1"
$wsWind.Range("M8").Value = 0.0015
$wsWind.Range("O8").Value = 31.372355844942916
$wsWind.Range("P8").Value = 1

# won-ITA_29 cost class 3, 4 & 5 (rows 14, 15, 16 rotate)
$wsWind.Range("C14").Value = "e_won-ITA_29_c5"
$wsWind.Range("D14").Value = "wind resource -- CF class won-ITA_29 -- cost class 5"
$wsWind.Range("K14").Value = "e_won-ITA_29_c5"
$wsWind.Range("P14").Value = 5

$wsWind.Range("C15").Value = "e_won-ITA_29_c3"
$wsWind.Range("D15").Value = "wind resource -- CF class won-ITA_29 -- cost class 3"
$wsWind.Range("K15").Value = "e_won-ITA_29_c3"
$wsWind.Range("P15").Value = 3

$wsWind.Range("C16").Value = "e_won-ITA_29_c4"
$wsWind.Range("D16").Value = "wind resource -- CF class won-ITA_29 -- cost class 4"
$wsWind.Range("K16").Value = "e_won-ITA_29_c4"
$wsWind.Range("P16").Value = 4

# won-ITA_25 cost class 3 & 4 (row 26 <-> row 27)
$wsWind.Range("C26").Value = "e_won-ITA_25_c4"
$wsWind.Range("D26").Value = "wind resource -- CF class won-ITA_25 -- cost class 4"
$wsWind.Range("K26").Value = "e_won-ITA_25_c4"
$wsWind.Range("M26").Value = 0.0045
$wsWind.Range("O26").Value = 107.03467623779731
$wsWind.Range("P26").Value = 4

$wsWind.Range("C27").Value = "e_won-ITA_25_c3"
$wsWind.Range("D27").Value = "wind resource -- CF class won-ITA_25 -- cost class 3"
$wsWind.Range("K27").Value = "e_won-ITA_25_c3"
$wsWind.Range("M27").Value = 3.2032500000000002
$wsWind.Range("O27").Value = 42.969340484621256
$wsWind.Range("P27").Value = 3

# won-ITA_11 cost class 2 & 5 (row 94 <-> row 95)
$wsWind.Range("C94").Value = "e_won-ITA_11_c2"
$wsWind.Range("D94").Value = "wind resource -- CF class won-ITA_11 -- cost class 2"
$wsWind.Range("K94").Value = "e_won-ITA_11_c2"
$wsWind.Range("M94").Value = 0.03
$wsWind.Range("O94").Value = 94.983314628615091
$wsWind.Range("P94").Value = 2

$wsWind.Range("C95").Value = "e_won-ITA_11_c5"
$wsWind.Range("D95").Value = "wind resource -- CF class won-ITA_11 -- cost class 5"
$wsWind.Range("K95").Value = "e_won-ITA_11_c5"
$wsWind.Range("M95").Value = 0.43575000000000003
$wsWind.Range("O95").Value = 120.30091898433189
$wsWind.Range("P95").Value = 5

# won-ITA_10 cost class 1 & 2 (row 98 <-> row 99)
$wsWind.Range("C98").Value = "e_won-ITA_10_c2"
$wsWind.Range("D98").Value = "wind resource -- CF class won-ITA_10 -- cost class 2"
$wsWind.Range("K98").Value = "e_won-ITA_10_c2"
$wsWind.Range("P98").Value = 2

$wsWind.Range("C99").Value = "e_won-ITA_10_c1"
$wsWind.Range("D99").Value = "wind resource -- CF class won-ITA_10 -- cost class 1"
$wsWind.Range("K99").Value = "e_won-ITA_10_c1"
$wsWind.Range("P99").Value = 1
